$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.748.99'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '1.878.75'
$ws.Range("E3").Value = '  +1.62%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.36%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4724'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3962'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.84'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08048'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.026'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.88'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.18%  '
$ws.Range("D13").Value = '1.900.75'
$ws.Range("E13").Value = '  +3.05%  '
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001050'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.94%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06619'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D22").Value = '27.773.68'
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("E24").Value = '  +2.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.301'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("D26").Value = '2.109.27'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.102'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.605'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9701'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09565'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.453'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.625'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.311'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.228'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.210'
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5996'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1906'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.259'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5699'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.940'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.68%  '
$ws.Range("E50").Value = '  +11.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06817'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.29%  '
